$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45954.01041666666
$ws.Range("A3").Value = 45954.02083333334
$ws.Range("A4").Value = 45954.03125
$ws.Range("A5").Value = 45954.04166666666
$ws.Range("A6").Value = 45954.05208333334
$ws.Range("A7").Value = 45954.0625
$ws.Range("A8").Value = 45954.07291666666
$ws.Range("A9").Value = 45954.08333333334
$ws.Range("A10").Value = 45954.09375
$ws.Range("A11").Value = 45954.10416666666
$ws.Range("A12").Value = 45954.11458333334
$ws.Range("A13").Value = 45954.125
$ws.Range("A14").Value = 45954.13541666666
$ws.Range("A15").Value = 45954.14583333334
$ws.Range("A16").Value = 45954.15625
$ws.Range("A17").Value = 45954.16666666666
$ws.Range("A18").Value = 45954.17708333334
$ws.Range("A19").Value = 45954.1875
$ws.Range("A20").Value = 45954.19791666666
$ws.Range("A21").Value = 45954.20833333334
$ws.Range("A22").Value = 45954.21875
$ws.Range("A23").Value = 45954.22916666666
$ws.Range("A24").Value = 45954.23958333334
$ws.Range("A25").Value = 45954.25
$ws.Range("A26").Value = 45954.26041666666
$ws.Range("A27").Value = 45954.27083333334
$ws.Range("A28").Value = 45954.28125
$ws.Range("A29").Value = 45954.29166666666
$ws.Range("A30").Value = 45954.30208333334
$ws.Range("A31").Value = 45954.3125
$ws.Range("A32").Value = 45954.32291666666
$ws.Range("A33").Value = 45954.33333333334
$ws.Range("A34").Value = 45954.34375
$ws.Range("A35").Value = 45954.35416666666
$ws.Range("A36").Value = 45954.36458333334
$ws.Range("A37").Value = 45954.375
$ws.Range("A38").Value = 45954.38541666666
$ws.Range("A39").Value = 45954.39583333334
$ws.Range("A40").Value = 45954.40625
$ws.Range("A41").Value = 45954.41666666666
$ws.Range("A42").Value = 45954.42708333334
$ws.Range("A43").Value = 45954.4375
$ws.Range("A44").Value = 45954.44791666666
$ws.Range("A45").Value = 45954.45833333334
$ws.Range("A46").Value = 45954.46875
$ws.Range("A47").Value = 45954.47916666666
$ws.Range("A48").Value = 45954.48958333334
$ws.Range("A49").Value = 45954.5
$ws.Range("A50").Value = 45954.51041666666
$ws.Range("A51").Value = 45954.52083333334
$ws.Range("A52").Value = 45954.53125
$ws.Range("A53").Value = 45954.54166666666
$ws.Range("A54").Value = 45954.55208333334
$ws.Range("A55").Value = 45954.5625
$ws.Range("A56").Value = 45954.57291666666
$ws.Range("A57").Value = 45954.58333333334
$ws.Range("A58").Value = 45954.59375
$ws.Range("A59").Value = 45954.60416666666
$ws.Range("A60").Value = 45954.61458333334
$ws.Range("A61").Value = 45954.625
$ws.Range("A62").Value = 45954.63541666666
$ws.Range("A63").Value = 45954.64583333334
$ws.Range("A64").Value = 45954.65625
$ws.Range("A65").Value = 45954.66666666666
$ws.Range("A66").Value = 45954.67708333334
$ws.Range("A67").Value = 45954.6875
$ws.Range("A68").Value = 45954.69791666666
$ws.Range("A69").Value = 45954.70833333334
$ws.Range("A70").Value = 45954.71875
$ws.Range("A71").Value = 45954.72916666666
$ws.Range("A72").Value = 45954.73958333334
$ws.Range("A73").Value = 45954.75
$ws.Range("A74").Value = 45954.76041666666
$ws.Range("A75").Value = 45954.77083333334
$ws.Range("A76").Value = 45954.78125
$ws.Range("A77").Value = 45954.79166666666
$ws.Range("A78").Value = 45954.80208333334
$ws.Range("A79").Value = 45954.8125
$ws.Range("A80").Value = 45954.82291666666
$ws.Range("A81").Value = 45954.83333333334
$ws.Range("A82").Value = 45954.84375
$ws.Range("A83").Value = 45954.85416666666
$ws.Range("A84").Value = 45954.86458333334
$ws.Range("A85").Value = 45954.875
$ws.Range("A86").Value = 45954.88541666666
$ws.Range("A87").Value = 45954.89583333334
$ws.Range("A88").Value = 45954.90625
$ws.Range("A89").Value = 45954.91666666666
$ws.Range("A90").Value = 45954.92708333334
$ws.Range("A91").Value = 45954.9375
$ws.Range("A92").Value = 45954.94791666666
$ws.Range("A93").Value = 45954.95833333334
$ws.Range("A94").Value = 45954.96875
$ws.Range("A95").Value = 45954.97916666666
$ws.Range("A96").Value = 45954.98958333334
$ws.Range("A97").Value = 45955.0

$ws.Range("B29").Value = 2
$ws.Range("B30").Value = 23
$ws.Range("B31").Value = 57
$ws.Range("B32").Value = 0
$ws.Range("B33").Value = 0
$ws.Range("B34").Value = 0
$ws.Range("B35").Value = 0
$ws.Range("B36").Value = 0
$ws.Range("B37").Value = 0
$ws.Range("B38").Value = 0
$ws.Range("B39").Value = 0
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("B42").Value = 0
